# Update the fs/txt change-method sample data in Sheet1:
#   column A now holds HTML file names instead of raw site URLs,
#   while columns B and C keep their existing dummy/test labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "test01.html"
$ws.Range("B1").Value = "dummy"
$ws.Range("C1").Value = "test01"

$ws.Range("A2").Value = "test02.html"
$ws.Range("B2").Value = "dummy"
$ws.Range("C2").Value = "test02"

# Move the active selection to B3, matching the saved workbook view.
[void]$ws.Range("B3").Select()
